$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A15").Value = "TFA"
$ws.Range("B15").Value = "Client Supervisório"
$ws.Range("C15").Value = "10.28.2.79"
$ws.Range("D15").Value = "Hostname"
$ws.Range("E15").Value = "./SCADA"
$ws.Range("F15").Value = "tfa@SS01"

$ws.Hyperlinks.Add($ws.Range("F15"), "mailto:tfa@SS01")

$ws.Range("F15").Select()
